# Updated cryptos list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.562.75"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'324.18"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4467"
$ws.Range("E7").Value = "  +5.01%  "
$ws.Range("D8").Value = "'0.3582"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").Value = "'0.07488"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'41.99"
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("D11").Value = "'1.094"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "'0.9999"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'20.83"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").Value = "'6.027"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").Value = "'7.116"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "'1.743.02"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "'93.20"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'0.06414"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "'0.9992"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'5.812"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").Value = "'27.609.63"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "'2.107"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").Value = "'162.77"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "'20.44"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "'1.952.90"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").Value = "'2.083"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("D30").Value = "'126.29"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "'1.079"
$ws.Range("E31").Value = "  -7.83%  "
$ws.Range("D32").Value = "'0.09082"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "'3.669"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("D34").Value = "'5.532"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").Value = "'11.97"
$ws.Range("E35").Value = "  -5.38%  "
$ws.Range("D36").Value = "'0.02289"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "'0.2096"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'0.06031"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "'0.6369"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "'4.958"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").Value = "'1.201"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").Value = "'1.376"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "'7.813"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Value = "'13.20"
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").Value = "'0.5911"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'3.711"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'122.30"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "'1.953"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'0.06854"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Value = "'72.44"
$ws.Range("E51").Value = "  -2.67%  "
